# Insert a new weekly price-report row at row 218 of the "Pepino ensalada"
# (Agro Chillan) sheet, pushing the existing rows 218-271 down to 219-272.
#
# Equivalent to: right-click row 218 header -> Insert, then fill in the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 218; everything below (incl. the
# old row 218) shifts down by one, and the sheet dimension grows to R272.
$ws.Rows.Item(218).Insert()

# Populate the newly-inserted row 218 with the new weekly observation.
$ws.Range("A218").Value = 7
$ws.Range("B218").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C218").Value = "Ñuble"
$ws.Range("D218").Value = 44889
$ws.Range("E218").Value = 16
$ws.Range("F218").Value = 100112043
$ws.Range("G218").Value = "Pepino ensalada"
$ws.Range("H218").Value = "Sin especificar"
$ws.Range("I218").Value = "Primera"
$ws.Range("J218").Value = 120
$ws.Range("K218").Value = 17000
$ws.Range("L218").Value = 18000
$ws.Range("M218").Value = 17500
$ws.Range("N218").Value = "$/caja 80 unidades"
$ws.Range("O218").Value = "Región del Maule"
$ws.Range("P218").Value = 219
$ws.Range("Q218").Value = 80
$ws.Range("R218").Value = "Hortaliza"
